$wb = $excel.ActiveWorkbook

# --- Fix mineral mass data entry error: values were recorded in mg, should be g ---
$bottle = $wb.Worksheets.Item("Bottle Results")
$bottle.Activate()

$massValues = @{
    2  = 43
    3  = 42
    4  = 38
    5  = 36
    6  = 42
    7  = 35
    8  = 36
    9  = 35
    10 = 37
    11 = 42
    12 = 48
    13 = 35
    14 = 34
    15 = 42
    16 = 33
    17 = 39
    18 = 49
    19 = 33
}

foreach ($row in $massValues.Keys) {
    $mg = $massValues[$row]
    $g = $mg / 1000
    $bottle.Range("E$row").Value = $g
}

# Restore the frozen-pane view / selection state for this sheet
$bottle.Range("E20").Select()

# --- Update view state (selection) on the other worksheets touched by this edit ---
$scint = $wb.Worksheets.Item("Scintillation Counter Results")
$scint.Activate()
$scint.Range("C77").Select()

$countActivity = $wb.Worksheets.Item("Count->Actual Activity")
$countActivity.Activate()
$countActivity.Range("F5").Select()

# --- Finish on Averaged Results (keeps it the active/selected tab, matching the workbook) ---
$averaged = $wb.Worksheets.Item("Averaged Results")
$averaged.Activate()
$averaged.Range("P15").Select()

# --- The calibration chart's Y axis was re-formatted from scientific notation to General ---
$calib = $wb.Worksheets.Item("Calibration Data")
$calib.Activate()
$chart = $calib.ChartObjects().Item(1).Chart
$yAxis = $chart.Axes(2)
$yAxis.TickLabels.NumberFormat = "General"

# Return focus to Averaged Results to match the workbook's final active tab/selection
$averaged.Activate()
$averaged.Range("P15").Select()
